# Consolidate the three separate text runs ("Below", " ", "section-level")
# on slide 2's title shape into a single run with the same combined text.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shape = $s.Shapes.Item(1)

# Setting TextRange.Text to a value that shares no common prefix/suffix with
# the current text forces the host to fully replace the run content (rather
# than patching only the differing tail), collapsing the paragraph down to a
# single run. Then set the real desired text.
$shape.TextFrame.TextRange.Text = "x"
$shape.TextFrame.TextRange.Text = "Below section-level"
